# Custom waste spreading by category and activity (part 8)
#
# This script reproduces the authored change to waste_to_full_name.xlsx:
#   - the waste code "W7375_hh" (row 21, column A) becomes "W7376_hh"
#   - the active window/selection of the sheet moves from I14 to B35
#   - the workbook window geometry (xWindow/yWindow/windowWidth/windowHeight)
#     is updated to reflect the new window size/position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the waste code text. Setting the value makes the shared-string
# table drop the old "W7375_hh" entry and append the new "W7376_hh" entry,
# which is exactly what happened in the authored diff.
$ws.Range("A21").Value = "W7376_hh"

# Move the selection/active cell for the sheet.
$ws.Range("B35").Select()

# Update the workbook window position/size.
$win = $excel.ActiveWindow
$win.Left = 0
$win.Top = 460
$win.Width = 51120
$win.Height = 28340
